# Prueba_Github.xlsx - "Add files via upload" edit
#
# Re-labels the content of both sheets' A1 cells and switches the active
# sheet/selection from Hoja1 (cell D4) to hoja2 (cell E4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("hoja2")

# New cell contents (was "hoja 1" / "hoja dos")
$ws1.Range("A1").Value = "datos hoja 1"
$ws2.Range("A1").Value = "datos hoja 2"

# Hoja1 no longer keeps the old D4 selection highlighted
$ws1.Range("A1").Select()

# hoja2 becomes the active tab, with E4 selected (was D4 on Hoja1)
$ws2.Activate()
$ws2.Range("E4").Select()
